# Adds two new columns (height, weight) to the right of the existing data
# and shifts the previously-computed "fantasy points" values into the new
# rightmost column, mirroring the behaviour of the upstream scraper re-run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New columns: height takes over E, weight takes over F, and the
# previous "fantasy points" header/values slide right into column G. -------
$ws.Range("E1").Value = "height"
$ws.Range("F1").Value = "weight"
$ws.Range("G1").Value = "fantasy points"

# Match the look of the existing bold/centered/bordered header cells
# (copy the formatting already used by the other header cells).
$ws.Range("D1").Copy()
$ws.Range("E1:G1").PasteSpecial(-4122)  # xlPasteFormats

# --- Data rows 2-7 ----------------------------------------------------------
# height/weight are constant for every row; the value that used to live in
# column E (fantasy points) moves over to the new column G.
$heightValue = 6.5
$weightValue = 257

$fantasyPoints = @{
    2 = 0
    3 = 6.5
    4 = 0
    5 = 3.9
    6 = 0.4
    7 = 0
}

foreach ($row in 2..7) {
    $ws.Range("E$row").Value = $heightValue
    $ws.Range("F$row").Value = $weightValue
    $ws.Range("G$row").Value = $fantasyPoints[$row]
}
